$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.348.01'
$ws.Range('E2').Value = '  -4.26%  '

$ws.Range('D3').Value = '2.485.17'
$ws.Range('E3').Value = '  -4.02%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.30'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.72%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.84'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.97%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('E8').Value = '  -3.45%  '

$ws.Range('D9').Value = '2.511.86'
$ws.Range('E9').Value = '  -3.26%  '

$ws.Range('E10').Value = '  -3.61%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.159'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.88%  '

$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.62'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.69%  '

$ws.Range('E13').Value = '  -2.15%  '

$ws.Range('D14').Value = '2.954.59'
$ws.Range('E14').Value = '  -3.06%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.17'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -5.70%  '

$ws.Range('D16').Value = '59.345.79'
$ws.Range('E16').Value = '  -4.16%  '

$ws.Range('E17').Value = '  -3.44%  '

$ws.Range('D18').Value = '2.514.85'
$ws.Range('E18').Value = '  -2.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.50'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.34%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.35'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.58'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.996'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.11%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.79'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.23%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.21'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.20%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.445'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -10.96%  '

$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.11%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.162'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.74%  '

$ws.Range('D28').Value = '2.604.18'
$ws.Range('E28').Value = '  -3.86%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.87'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.14'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.82%  '

$ws.Range('E31').Value = '  -1.31%  '

$ws.Range('E32').Value = '  -6.47%  '

$ws.Range('E33').Value = '  -5.08%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.51'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.94%  '

$ws.Range('E36').Value = '  -0.29%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.62'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.41%  '

$ws.Range('E38').Value = '  -5.16%  '

$ws.Range('E39').Value = '  -6.42%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.92'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.90%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '314.88'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -7.54%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.80'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.92%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.75'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.833'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.17%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.995'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.28%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.598'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.55%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.76'
$ws.Range('D47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.99'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.74%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0527'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.01%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0932'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.57%  '

$ws.Range('E51').Value = '  -4.17%  '
